$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$t = $ws.ListObjects.Item(1)
$rows = $t.ListRows()
$newRow = $rows.Add()

$ws.Range("A84:J84").Copy()
$ws.Range("A85:J85").PasteSpecial(-4122)

$ws.Range("A85").Value = 43985
$ws.Range("B85").Value = 81333
$ws.Range("C85").Value = 828
$ws.Range("D85").Value = 1477
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 5
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 109
$ws.Range("J85").Value = 0

$ws.Range("A85:J85").Select()
